# Update the author byline on slide 1 to add the new co-author
# "Gabriel de Souza Franco" (previously just "Gabriel").
#
# Slide 1, shape 2 (the subtitle placeholder) holds two paragraphs:
#   "January 11, 2020"
#   "Diane Scherpereel, Gabriel, Roger J. Albarran"
# The second paragraph is made up of four runs:
#   R1 "Diane "
#   R2 "Scherpereel"
#   R3 ", Gabriel, Roger J. "
#   R4 "Albarran"
#
# Target wording: "Diane Scherpereel, Gabriel de Souza Franco, Roger J. Albarran"
# achieved by keeping the same four runs but changing their text to:
#   R1 "Diane Scherpereel, Gabriel "
#   R2 "de Souza Franco"
#   R3 ", Roger J. "
#   R4 "Albarran"                 (unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$full = $tr.Text

$oldFirst  = "Diane "
$oldName   = "Scherpereel"
$oldMiddle = ", Gabriel, Roger J. "

# 1-based TextRange.Characters start positions of each run, located dynamically
# from the current text so the script isn't tied to hard-coded offsets.
$startFirst  = $full.IndexOf($oldFirst) + 1
$startName   = $full.IndexOf($oldName) + 1
$startMiddle = $full.IndexOf($oldMiddle) + 1

# Apply the edits right-to-left so offsets computed above (which refer to the
# original text) stay valid while earlier runs are still being resized.
$tr.Characters($startMiddle, $oldMiddle.Length).Text = ", Roger J. "
$tr.Characters($startName, $oldName.Length).Text = "de Souza Franco"
$tr.Characters($startFirst, $oldFirst.Length).Text = "Diane Scherpereel, Gabriel "
